$d = $word.ActiveDocument

# --- Edit paragraph 1: "Requirements for phiTsUtras" -> split run + proofErr spell markers ---
$xmlPara1 = @'
<w:body xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:p>
      <w:r>
        <w:t xml:space="preserve">Requirements for </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>phiTsUtras</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
    </w:p>
</w:body>
'@
$p1 = $d.Paragraphs(1).Range
$p1.InsertXML($xmlPara1)

# --- Replace paragraphs 3..end with the fully updated/expanded content ---
$xmlRest = @'
<w:body xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:p>
      <w:pPr>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="11"/>
        </w:numPr>
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">Ends with </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>dakziRa</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
    </w:p>
<w:p>
      <w:pPr>
        <w:numPr>
          <w:ilvl w:val="1"/>
          <w:numId w:val="11"/>
        </w:numPr>
      </w:pPr>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Mangal" w:hAnsi="Mangal" w:cs="Mangal"/>
        </w:rPr>
        <w:t>साधुत्ववाचि</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
    </w:p>
<w:p>
      <w:pPr>
        <w:numPr>
          <w:ilvl w:val="1"/>
          <w:numId w:val="11"/>
        </w:numPr>
      </w:pPr>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Mangal" w:hAnsi="Mangal" w:cs="Mangal"/>
        </w:rPr>
        <w:t>स्वाङ्गवाचि</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
    </w:p>
<w:p>
      <w:pPr>
        <w:numPr>
          <w:ilvl w:val="1"/>
          <w:numId w:val="11"/>
        </w:numPr>
      </w:pPr>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Mangal" w:hAnsi="Mangal" w:cs="Mangal"/>
        </w:rPr>
        <w:t>अन्य</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
    </w:p>
<w:p>
      <w:pPr>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="11"/>
        </w:numPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Mangal" w:hAnsi="Mangal" w:cs="Mangal"/>
        </w:rPr>
        <w:t xml:space="preserve">Ends with </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Mangal" w:hAnsi="Mangal" w:cs="Mangal"/>
        </w:rPr>
        <w:t>kfzRa</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
    </w:p>
<w:p>
      <w:pPr>
        <w:numPr>
          <w:ilvl w:val="1"/>
          <w:numId w:val="11"/>
        </w:numPr>
      </w:pPr>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Mangal" w:hAnsi="Mangal" w:cs="Mangal"/>
        </w:rPr>
        <w:t>मृगाख्या</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
    </w:p>
<w:p>
      <w:pPr>
        <w:numPr>
          <w:ilvl w:val="1"/>
          <w:numId w:val="11"/>
        </w:numPr>
      </w:pPr>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Mangal" w:hAnsi="Mangal" w:cs="Mangal"/>
        </w:rPr>
        <w:t>नाम</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
    </w:p>
<w:p>
      <w:pPr>
        <w:numPr>
          <w:ilvl w:val="1"/>
          <w:numId w:val="11"/>
        </w:numPr>
      </w:pPr>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Mangal" w:hAnsi="Mangal" w:cs="Mangal"/>
        </w:rPr>
        <w:t>अन्य</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
    </w:p>
<w:p>
      <w:pPr>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="11"/>
        </w:numPr>
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">Ends with </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:t>arjuna</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:t xml:space="preserve"> // This is avoidable, because the meaning of </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>tRNa</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> is only seen in neuter gender. </w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:t>Rest are</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:t xml:space="preserve"> not </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>tRNAkhyA</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve">. Let us check the lexicon whether we get some other meaning of </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:t>arjuna</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:t xml:space="preserve"> neuter. If not, this can be removed.</w:t>
      </w:r>
    </w:p>
<w:p>
      <w:pPr>
        <w:numPr>
          <w:ilvl w:val="1"/>
          <w:numId w:val="11"/>
        </w:numPr>
      </w:pPr>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Mangal" w:hAnsi="Mangal" w:cs="Mangal"/>
        </w:rPr>
        <w:t>तृणाख्या</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
    </w:p>
<w:p>
      <w:pPr>
        <w:numPr>
          <w:ilvl w:val="1"/>
          <w:numId w:val="11"/>
        </w:numPr>
      </w:pPr>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Mangal" w:hAnsi="Mangal" w:cs="Mangal"/>
        </w:rPr>
        <w:t>अन्य</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
    </w:p>
<w:p>
      <w:pPr>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="11"/>
        </w:numPr>
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">Ends with </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>Arya</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
    </w:p>
<w:p>
      <w:pPr>
        <w:numPr>
          <w:ilvl w:val="1"/>
          <w:numId w:val="11"/>
        </w:numPr>
      </w:pPr>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Mangal" w:hAnsi="Mangal" w:cs="Mangal"/>
        </w:rPr>
        <w:t>स्वाम्याख्या</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
    </w:p>
<w:p>
      <w:pPr>
        <w:numPr>
          <w:ilvl w:val="1"/>
          <w:numId w:val="11"/>
        </w:numPr>
      </w:pPr>
      <w:bookmarkStart w:id="0" w:name="_GoBack"/>
      <w:bookmarkEnd w:id="0"/>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Mangal" w:hAnsi="Mangal" w:cs="Mangal"/>
        </w:rPr>
        <w:lastRenderedPageBreak/>
        <w:t>अन्य</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
    </w:p>
<w:p>
      <w:pPr>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="11"/>
        </w:numPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Mangal" w:hAnsi="Mangal" w:cs="Mangal"/>
        </w:rPr>
        <w:t xml:space="preserve">Ends with </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Mangal" w:hAnsi="Mangal" w:cs="Mangal"/>
        </w:rPr>
        <w:t>ASA</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
    </w:p>
<w:p>
      <w:pPr>
        <w:numPr>
          <w:ilvl w:val="1"/>
          <w:numId w:val="11"/>
        </w:numPr>
      </w:pPr>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Mangal" w:hAnsi="Mangal" w:cs="Mangal"/>
        </w:rPr>
        <w:t>दिगाख्या</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
    </w:p>
<w:p>
      <w:pPr>
        <w:numPr>
          <w:ilvl w:val="1"/>
          <w:numId w:val="11"/>
        </w:numPr>
      </w:pPr>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Mangal" w:hAnsi="Mangal" w:cs="Mangal"/>
        </w:rPr>
        <w:t>अदिगाख्या</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
    </w:p>
<w:p>
      <w:pPr>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="11"/>
        </w:numPr>
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">Ends with </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>jyezWa</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve">, </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>kanizWa</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
    </w:p>
<w:p>
      <w:pPr>
        <w:numPr>
          <w:ilvl w:val="1"/>
          <w:numId w:val="11"/>
        </w:numPr>
      </w:pPr>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Mangal" w:hAnsi="Mangal" w:cs="Mangal"/>
        </w:rPr>
        <w:t>वयसि</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
    </w:p>
<w:p>
      <w:pPr>
        <w:numPr>
          <w:ilvl w:val="1"/>
          <w:numId w:val="11"/>
        </w:numPr>
      </w:pPr>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Mangal" w:hAnsi="Mangal" w:cs="Mangal"/>
        </w:rPr>
        <w:t>अन्य</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
    </w:p>
<w:p>
      <w:pPr>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="11"/>
        </w:numPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Mangal" w:hAnsi="Mangal" w:cs="Mangal"/>
        </w:rPr>
        <w:t xml:space="preserve">has </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Mangal" w:hAnsi="Mangal" w:cs="Mangal"/>
        </w:rPr>
        <w:t>arr</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Mangal" w:hAnsi="Mangal" w:cs="Mangal"/>
        </w:rPr>
        <w:t>($text,'/[</w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Mangal" w:hAnsi="Mangal" w:cs="Mangal"/>
        </w:rPr>
        <w:t>uU</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Mangal" w:hAnsi="Mangal" w:cs="Mangal"/>
        </w:rPr>
        <w:t>][</w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Mangal" w:hAnsi="Mangal" w:cs="Mangal"/>
        </w:rPr>
        <w:t>KPCWTcwtkp</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Mangal" w:hAnsi="Mangal" w:cs="Mangal"/>
        </w:rPr>
        <w:t>]/')</w:t>
      </w:r>
    </w:p>
<w:p>
      <w:pPr>
        <w:numPr>
          <w:ilvl w:val="1"/>
          <w:numId w:val="11"/>
        </w:numPr>
      </w:pPr>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Mangal" w:hAnsi="Mangal" w:cs="Mangal"/>
        </w:rPr>
        <w:t>कृत्रिमाख्या</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
    </w:p>
<w:p>
      <w:pPr>
        <w:numPr>
          <w:ilvl w:val="1"/>
          <w:numId w:val="11"/>
        </w:numPr>
      </w:pPr>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Mangal" w:hAnsi="Mangal" w:cs="Mangal"/>
        </w:rPr>
        <w:t>अन्य</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
    </w:p>
<w:p>
      <w:pPr>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="11"/>
        </w:numPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Mangal" w:hAnsi="Mangal" w:cs="Mangal"/>
        </w:rPr>
        <w:t xml:space="preserve">first vowel is a </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Mangal" w:hAnsi="Mangal" w:cs="Mangal"/>
        </w:rPr>
        <w:t>hrasva</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Mangal" w:hAnsi="Mangal" w:cs="Mangal"/>
        </w:rPr>
        <w:t xml:space="preserve"> vowel except ‘f’ and the last vowel is </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Mangal" w:hAnsi="Mangal" w:cs="Mangal"/>
        </w:rPr>
        <w:t>hrasva</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
    </w:p>
<w:p>
      <w:pPr>
        <w:numPr>
          <w:ilvl w:val="1"/>
          <w:numId w:val="11"/>
        </w:numPr>
      </w:pPr>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Mangal" w:hAnsi="Mangal" w:cs="Mangal"/>
        </w:rPr>
        <w:t>ताच्छील्ये</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
    </w:p>
<w:p>
      <w:pPr>
        <w:numPr>
          <w:ilvl w:val="1"/>
          <w:numId w:val="11"/>
        </w:numPr>
      </w:pPr>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Mangal" w:hAnsi="Mangal" w:cs="Mangal"/>
        </w:rPr>
        <w:t>अन्य</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
    </w:p>
<w:p>
      <w:pPr>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="11"/>
        </w:numPr>
      </w:pPr>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Mangal" w:hAnsi="Mangal" w:cs="Mangal"/>
        </w:rPr>
        <w:t>akza</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
    </w:p>
<w:p>
      <w:pPr>
        <w:numPr>
          <w:ilvl w:val="1"/>
          <w:numId w:val="11"/>
        </w:numPr>
      </w:pPr>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Mangal" w:hAnsi="Mangal" w:cs="Mangal"/>
        </w:rPr>
        <w:t>देवन</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
    </w:p>
<w:p>
      <w:pPr>
        <w:numPr>
          <w:ilvl w:val="1"/>
          <w:numId w:val="11"/>
        </w:numPr>
      </w:pPr>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Mangal" w:hAnsi="Mangal" w:cs="Mangal"/>
        </w:rPr>
        <w:lastRenderedPageBreak/>
        <w:t>अदेवने</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
    </w:p>
<w:p>
      <w:pPr>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="11"/>
        </w:numPr>
      </w:pPr>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>arDa</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
    </w:p>
<w:p>
      <w:pPr>
        <w:numPr>
          <w:ilvl w:val="1"/>
          <w:numId w:val="11"/>
        </w:numPr>
      </w:pPr>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Mangal" w:hAnsi="Mangal" w:cs="Mangal"/>
        </w:rPr>
        <w:t>असमद्योतने</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
    </w:p>
<w:p>
      <w:pPr>
        <w:numPr>
          <w:ilvl w:val="1"/>
          <w:numId w:val="11"/>
        </w:numPr>
      </w:pPr>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Mangal" w:hAnsi="Mangal" w:cs="Mangal"/>
        </w:rPr>
        <w:t>समद्योतने</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
    </w:p>
<w:p>
      <w:pPr>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="11"/>
        </w:numPr>
      </w:pPr>
      <w:r>
        <w:t>Ends with [</w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>Dy</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t>]A</w:t>
      </w:r>
    </w:p>
<w:p>
      <w:pPr>
        <w:numPr>
          <w:ilvl w:val="1"/>
          <w:numId w:val="11"/>
        </w:numPr>
      </w:pPr>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Mangal" w:hAnsi="Mangal" w:cs="Mangal"/>
        </w:rPr>
        <w:t>नित्यस्त्रीलिङ्ग</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
    </w:p>
<w:p>
      <w:pPr>
        <w:numPr>
          <w:ilvl w:val="1"/>
          <w:numId w:val="11"/>
        </w:numPr>
      </w:pPr>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Mangal" w:hAnsi="Mangal" w:cs="Mangal"/>
        </w:rPr>
        <w:t>अन्य</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
    </w:p>
<w:p>
      <w:pPr>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="11"/>
        </w:numPr>
      </w:pPr>
      <w:r>
        <w:t>Ends with ‘A’ / ‘I’ and doesn’t end with [</w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>Dy</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t>]A</w:t>
      </w:r>
    </w:p>
<w:p>
      <w:pPr>
        <w:numPr>
          <w:ilvl w:val="1"/>
          <w:numId w:val="11"/>
        </w:numPr>
      </w:pPr>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Mangal" w:hAnsi="Mangal" w:cs="Mangal"/>
        </w:rPr>
        <w:t>नित्यस्त्रीलिङ्ग</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
    </w:p>
<w:p>
      <w:pPr>
        <w:numPr>
          <w:ilvl w:val="1"/>
          <w:numId w:val="11"/>
        </w:numPr>
      </w:pPr>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Mangal" w:hAnsi="Mangal" w:cs="Mangal"/>
        </w:rPr>
        <w:t>अन्य</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
    </w:p>
<w:p>
      <w:pPr>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="11"/>
        </w:numPr>
      </w:pPr>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Mangal" w:hAnsi="Mangal" w:cs="Mangal"/>
        </w:rPr>
        <w:t>Sukla</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Mangal" w:hAnsi="Mangal" w:cs="Mangal"/>
        </w:rPr>
        <w:t xml:space="preserve"> / </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Mangal" w:hAnsi="Mangal" w:cs="Mangal"/>
        </w:rPr>
        <w:t>gOra</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
    </w:p>
<w:p>
      <w:pPr>
        <w:numPr>
          <w:ilvl w:val="1"/>
          <w:numId w:val="11"/>
        </w:numPr>
      </w:pPr>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Mangal" w:hAnsi="Mangal" w:cs="Mangal"/>
        </w:rPr>
        <w:t>नामवाचि</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
    </w:p>
<w:p>
      <w:pPr>
        <w:numPr>
          <w:ilvl w:val="1"/>
          <w:numId w:val="11"/>
        </w:numPr>
      </w:pPr>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Mangal" w:hAnsi="Mangal" w:cs="Mangal"/>
        </w:rPr>
        <w:t>अन्य</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
    </w:p>
<w:p>
      <w:pPr>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="11"/>
        </w:numPr>
      </w:pPr>
    </w:p>
</w:body>
'@
$startRange = $d.Paragraphs(3).Range.Start
$lastIndex = $d.Paragraphs.Count
$endRange = $d.Paragraphs($lastIndex).Range.End
$r = $d.Range($startRange, $endRange)
$r.InsertXML($xmlRest)
